# Fruta / hortaliza, semanal
# Weekly refresh of the price data: rows 4-9 (columns D, J, K, L, M, N, O, P, Q)
# are re-shuffled to reflect the updated weekly snapshot. Columns A, B, C, E,
# F, G, H, I, R are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for rows 4-9, keyed by original row number.
$rows = @{
    4 = @{ D = 44215; J = 140; K = 3500; L = 4000; M = 3768; N = "`$/paquete 2 kilos";     O = "Provincia de Diguillín"; P = 1884; Q = 2 }
    5 = @{ D = 44166; J = 70;  K = 3500; L = 4000; M = 3679; N = "`$/paquete 36 unidades"; O = "Región Metropolitana";    P = 102;  Q = 36 }
    6 = @{ D = 44209; J = 150; K = 3500; L = 4000; M = 3767; N = "`$/paquete 2 kilos";     O = "Provincia de Diguillín"; P = 1884; Q = 2 }
    7 = @{ D = 44210; J = 105; K = 3500; L = 4000; M = 3714; N = "`$/paquete 2 kilos";     O = "Provincia de Diguillín"; P = 1857; Q = 2 }
    8 = @{ D = 44161; J = 50;  K = 2800; L = 3000; M = 2900; N = "`$/paquete 2 kilos";     O = "Provincia de Diguillín"; P = 1450; Q = 2 }
    9 = @{ D = 44160; J = 43;  K = 3500; L = 4000; M = 3709; N = "`$/paquete 36 unidades"; O = "Región Metropolitana";    P = 103;  Q = 36 }
}

# Mapping of new row number -> source row number (i.e. new row 4 gets the
# values that used to live in row 5, etc).
$mapping = @{
    4 = 5
    5 = 4
    6 = 8
    7 = 9
    8 = 7
    9 = 6
}

foreach ($destRow in 4..9) {
    $srcRow = $mapping[$destRow]
    $vals = $rows[$srcRow]

    $ws.Range("D$destRow").Value2 = $vals.D
    $ws.Range("J$destRow").Value2 = $vals.J
    $ws.Range("K$destRow").Value2 = $vals.K
    $ws.Range("L$destRow").Value2 = $vals.L
    $ws.Range("M$destRow").Value2 = $vals.M
    $ws.Range("N$destRow").Value2 = $vals.N
    $ws.Range("O$destRow").Value2 = $vals.O
    $ws.Range("P$destRow").Value2 = $vals.P
    $ws.Range("Q$destRow").Value2 = $vals.Q
}
